$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date header in column AM (column 39), continuing the series
# of day headers that already populate columns C..AL.
$ws.Cells.Item(1, 39).Value = "27-jul"

# Fill in the new day's values for each data row (2-18), mirroring the
# values that were appended to the source file for this new "27-jul" column.
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(3, 39).Value = 16.408288513876712
$ws.Cells.Item(4, 39).Value = 14.803082072462523
$ws.Cells.Item(5, 39).Value = 29.697120882378265
$ws.Cells.Item(6, 39).Value = 0
$ws.Cells.Item(7, 39).Value = 5.5463615989111208
$ws.Cells.Item(8, 39).Value = 4.872284924428306
$ws.Cells.Item(9, 39).Value = 19.252245655744399
$ws.Cells.Item(10, 39).Value = 19.810990297932474
$ws.Cells.Item(11, 39).Value = 13.083454944338438
$ws.Cells.Item(12, 39).Value = 0
$ws.Cells.Item(13, 39).Value = 14.107237455014955
$ws.Cells.Item(14, 39).Value = 0
$ws.Cells.Item(15, 39).Value = 0
$ws.Cells.Item(16, 39).Value = 3.9678216052448567
$ws.Cells.Item(17, 39).Value = 0
$ws.Cells.Item(18, 39).Value = 0

# Match the formatting used by the other date-header cells (text format)
# for the new header cell in row 1.
$ws.Cells.Item(1, 39).NumberFormat = $ws.Cells.Item(1, 38).NumberFormat

# Move/restore the active selection, as happened in the source edit.
$ws.Range("AO7").Select()
